$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each existing data row (12-18) shifts its values down by one row (row N gets
# the values that used to live in row N-1), and the values that used to be in
# row 18 become a brand-new row 19. Row 12 receives a new data point that
# didn't exist before.

# Row 12 - new data point
$ws.Range("D12").Value = 44427
$ws.Range("J12").Value = 250
$ws.Range("K12").Value = 1300
$ws.Range("M12").Value = 1400
$ws.Range("P12").Value = 700

# Row 13 - previously row 12's values
$ws.Range("D13").Value = 44257
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 1400
$ws.Range("L13").Value = 1500
$ws.Range("M13").Value = 1450
$ws.Range("P13").Value = 725

# Row 14 - previously row 13's values
$ws.Range("D14").Value = 44390
$ws.Range("J14").Value = 250
$ws.Range("K14").Value = 2400
$ws.Range("L14").Value = 2500
$ws.Range("M14").Value = 2450
$ws.Range("P14").Value = 1225

# Row 15 - previously row 14's values
$ws.Range("D15").Value = 44172
$ws.Range("J15").Value = 200
$ws.Range("K15").Value = 1300
$ws.Range("L15").Value = 1500
$ws.Range("M15").Value = 1400
$ws.Range("P15").Value = 700

# Row 16 - previously row 15's values
$ws.Range("D16").Value = 44202
$ws.Range("J16").Value = 250
$ws.Range("K16").Value = 1800
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = 1900
$ws.Range("N16").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("P16").Value = 950
$ws.Range("Q16").Value = 2

# Row 17 - previously row 16's values
$ws.Range("D17").Value = 44181
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 1000
$ws.Range("L17").Value = 1200
$ws.Range("M17").Value = 1100
$ws.Range("N17").Value = "$/atado"
$ws.Range("P17").Value = 1100
$ws.Range("Q17").Value = 1

# Row 18 - previously row 17's values
$ws.Range("D18").Value = 44161
$ws.Range("J18").Value = 270
$ws.Range("K18").Value = 900
$ws.Range("L18").Value = 1000
$ws.Range("M18").Value = 950
$ws.Range("N18").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("P18").Value = 475
$ws.Range("Q18").Value = 2

# Row 19 - brand new row, holding what used to be row 18's values
$ws.Range("A19").Value = 1
$ws.Range("B19").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C19").Value = "Arica y Parinacota"
$ws.Range("D19").Value = 44253
$ws.Range("D19").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E19").Value = 15
$ws.Range("F19").Value = 100112044
$ws.Range("G19").Value = "Perejil"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 250
$ws.Range("K19").Value = 1800
$ws.Range("L19").Value = 2000
$ws.Range("M19").Value = 1900
$ws.Range("N19").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O19").Value = "Región de Arica y Parinacota"
$ws.Range("P19").Value = 950
$ws.Range("Q19").Value = 2
$ws.Range("R19").Value = "Hortaliza"
